$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text storage (avoid Excel auto-converting numeric-looking strings)
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @(
  @{ Row = 2; D = "26.213.55"; E = "  -0.27%  " },
  @{ Row = 3; D = "1.659.02"; E = "  -0.55%  " },
  @{ Row = 4; D = "1.004"; E = "  -0.39%  " },
  @{ Row = 5; D = "219.51"; E = "  -0.17%  " },
  @{ Row = 6; D = "0.5273"; E = "  -0.32%  " },
  @{ Row = 7; E = "  -0.33%  " },
  @{ Row = 9; D = "0.06387"; E = "  +0.29%  " },
  @{ Row = 10; D = "20.66"; E = "  -1.39%  " },
  @{ Row = 11; D = "0.07690"; E = "  -1.78%  " },
  @{ Row = 12; D = "4.624"; E = "  +2.21%  " },
  @{ Row = 13; D = "1.668.92"; E = "  -0.60%  " },
  @{ Row = 14; D = "1.888.15"; E = "  -0.43%  " },
  @{ Row = 15; D = "0.5647"; E = "  +0.95%  " },
  @{ Row = 16; D = ("0.0{0}8268" -f [char]8325); E = "  +2.09%  " },
  @{ Row = 17; D = "65.86"; E = "  +0.18%  " },
  @{ Row = 18; D = "26.201.18" },
  @{ Row = 19; D = "1.004"; E = "  -0.40%  " },
  @{ Row = 20; D = "4.692"; E = "  -0.62%  " },
  @{ Row = 21; D = "10.42"; E = "  +1.46%  " },
  @{ Row = 22; D = "192.00"; E = "  -4.18%  " },
  @{ Row = 23; D = "6.005"; E = "  -0.92%  " },
  @{ Row = 24; E = "  -0.36%  " },
  @{ Row = 25; D = "146.05"; E = "  -0.24%  " },
  @{ Row = 27; E = "  +1.03%  " },
  @{ Row = 28; E = "  -0.68%  " },
  @{ Row = 29; D = "1.526"; E = "  -0.15%  " },
  @{ Row = 30; D = "0.05658"; E = "  -3.75%  " },
  @{ Row = 31; D = "1.281"; E = "  -0.17%  " },
  @{ Row = 32; D = "3.497"; E = "  -0.31%  " },
  @{ Row = 33; D = "3.405"; E = "  +2.16%  " },
  @{ Row = 34; E = "  -0.89%  " },
  @{ Row = 35; D = "0.9534"; E = "  -1.08%  " },
  @{ Row = 36; D = "2.793"; E = "  -1.00%  " },
  @{ Row = 37; D = "2.405"; E = "  -0.98%  " },
  @{ Row = 38; D = "0.5780"; E = "  -0.39%  " },
  @{ Row = 40; D = "5.986"; E = "  +0.14%  " },
  @{ Row = 41; E = "  -0.33%  " },
  @{ Row = 42; D = "0.8371"; E = "  -2.39%  " },
  @{ Row = 43; D = "1.031.94"; E = "  -4.39%  " },
  @{ Row = 44; E = "  -1.23%  " },
  @{ Row = 45; D = "1.798.76"; E = "  -0.52%  " },
  @{ Row = 46; D = "58.61"; E = "  +0.00%  " },
  @{ Row = 47; E = "  +3.91%  " },
  @{ Row = 48; E = "  -0.82%  " },
  @{ Row = 49; D = "0.05345"; E = "  +3.92%  " },
  @{ Row = 50; D = "8.094"; E = "  +0.63%  " },
  @{ Row = 51; D = "0.4346"; E = "  -1.58%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}

# Restore default style/number format (no explicit style attribute, like the original cells)
$ws.Range("D2:E51").Style = "Normal"
